$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two percentage values
$ws.Range("B4").Value = 1.9
$ws.Range("B5").Value = 0.39

# Update the registration date text (cell is formatted as Text, so it stays a string)
$ws.Range("B7").Value = "01.01.2024"

# Update selection to match the new active range
$ws.Range("B4:B7").Select()
